$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their original text formatting
# so numeric-looking strings (e.g. "1.003", "23.752.66") and percent strings
# with surrounding spaces are not re-interpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '23.734.47'
$ws.Range('E2').Value = '  +1.47%  '
$ws.Range('D3').Value = '1.656.15'
$ws.Range('E3').Value = '  +1.58%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Value = '1.002'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Value = '303.45'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').Value = '0.3800'
$ws.Range('E7').Value = '  +0.85%  '
$ws.Range('D9').Value = '51.00'
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('E10').Value = '  +2.75%  '
$ws.Range('D11').Value = '0.08237'
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = '22.67'
$ws.Range('E13').Value = '  +1.95%  '
$ws.Range('D14').Value = '6.545'
$ws.Range('E14').Value = '  +1.15%  '
$ws.Range('E15').Value = '  +2.22%  '
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = '1.657.97'
$ws.Range('E17').Value = '  +2.33%  '
$ws.Range('D18').Value = '97.81'
$ws.Range('E18').Value = '  +3.24%  '
$ws.Range('D19').Value = '0.06977'
$ws.Range('E19').Value = '  +0.47%  '
$ws.Range('D20').Value = '6.812'
$ws.Range('E20').Value = '  +4.20%  '
$ws.Range('D21').Value = '17.77'
$ws.Range('E21').Value = '  +1.22%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('E23').Value = '  +2.59%  '
$ws.Range('D24').Value = '23.788.13'
$ws.Range('E24').Value = '  +1.75%  '
$ws.Range('D25').Value = '2.549'
$ws.Range('E25').Value = '  +1.53%  '
$ws.Range('D26').Value = '3.069'
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').Value = '21.34'
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('D28').Value = '151.74'
$ws.Range('E28').Value = '  +0.84%  '
$ws.Range('D29').Value = '5.217'
$ws.Range('E29').Value = '  -0.91%  '
$ws.Range('E30').Value = '  +1.26%  '
$ws.Range('D31').Value = '1.838.93'
$ws.Range('E31').Value = '  +2.21%  '
$ws.Range('D32').Value = '6.929'
$ws.Range('E32').Value = '  +4.71%  '
$ws.Range('D33').Value = '2.183'
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('E34').Value = '  +1.91%  '
$ws.Range('D35').Value = '11.81'
$ws.Range('E35').Value = '  +4.15%  '
$ws.Range('D36').Value = '0.02825'
$ws.Range('E36').Value = '  +2.11%  '
$ws.Range('E38').Value = '  +2.58%  '
$ws.Range('D39').Value = '0.08819'
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('D40').Value = '0.07100'
$ws.Range('E40').Value = '  -0.48%  '
$ws.Range('D41').Value = '13.36'
$ws.Range('E41').Value = '  +11.35%  '
$ws.Range('D42').Value = '0.7087'
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('D43').Value = '1.343'
$ws.Range('E43').Value = '  +0.94%  '
$ws.Range('D44').Value = '15.98'
$ws.Range('E44').Value = '  +0.70%  '
$ws.Range('D45').Value = '0.6551'
$ws.Range('E45').Value = '  +1.36%  '
$ws.Range('D46').Value = '2.335'
$ws.Range('E46').Value = '  +2.60%  '
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').Value = '3.965'
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('D49').Value = '0.07958'
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('D50').Value = '128.17'
$ws.Range('E50').Value = '  +1.39%  '
$ws.Range('E51').Value = '  +0.56%  '
